# Update the LR-pairs (Thbs1-Itga6) sheet with the new TPM-derived values.
# This mirrors the commit "update scripts wuth new tpm": the ligand/receptor
# average & total expression values were recomputed, which cascades into the
# derived-specificity and edge-weight columns as well.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row -> column -> new value, as found in the target workbook.
$updates = @{
    2  = @{ G = 3.075165666666667; H = 9.225497000000001; I = 0.02641273658732285; J = 0.02641273658732285;
            M = 139.728498; N = 419.185494; O = 0.9065295391216045; P = 0.9065295391216045;
            Q = 429.6882797045021; R = 3867.194517340518; S = 0.02394392592544612; T = 0.02394392592544612 }
    3  = @{ G = 3.075165666666667; H = 9.225497000000001; I = 0.02641273658732285; J = 0.02641273658732285;
            O = 0.005362677585431591; P = 0.005362677585431591;
            Q = 2.541869411698112; R = 22.87682470528301; S = 0.0001416429904667451; T = 0.0001416429904667451 }
    4  = @{ G = 3.075165666666667; H = 9.225497000000001; I = 0.02641273658732285; J = 0.02641273658732285;
            O = 0.0881077832929639; P = 0.0881077832929639;
            Q = 41.76243596917389; R = 375.861923722565; S = 0.00232716767140998; T = 0.00232716767140998 }
    5  = @{ I = 0.549422396165273; J = 0.5494223961652731;
            M = 139.728498; N = 419.185494; O = 0.9065295391216045; P = 0.9065295391216045;
            Q = 8938.125871920882; R = 80443.13284728794; S = 0.4980676315787925; T = 0.4980676315787926 }
    6  = @{ I = 0.549422396165273; J = 0.5494223961652731;
            O = 0.005362677585431591; P = 0.005362677585431591;
            S = 0.002946375168849625; T = 0.002946375168849626 }
    7  = @{ I = 0.549422396165273; J = 0.5494223961652731;
            O = 0.0881077832929639; P = 0.0881077832929639;
            S = 0.04840838941763084; T = 0.04840838941763084 }
    8  = @{ I = 0.424164867247404; J = 0.4241648672474041;
            M = 139.728498; N = 419.185494; O = 0.9065295391216045; P = 0.9065295391216045;
            Q = 6900.40850239286; R = 62103.67652153574; S = 0.3845179816173657; T = 0.3845179816173658 }
    9  = @{ I = 0.424164867247404; J = 0.4241648672474041;
            O = 0.005362677585431591; P = 0.005362677585431591;
            S = 0.00227465942611522; T = 0.00227465942611522 }
    10 = @{ I = 0.424164867247404; J = 0.4241648672474041;
            O = 0.0881077832929639; P = 0.0881077832929639;
            S = 0.03737222620392307; T = 0.03737222620392308 }
}

foreach ($row in $updates.Keys) {
    $cols = $updates[$row]
    foreach ($col in $cols.Keys) {
        $ws.Range("$col$row").Value = $cols[$col]
    }
}
